$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "sub_003"
$ws.Range("B4").Value = $false

$ws.Range("B3").Select()
